$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Price (column D) updates: force text storage so numeric-looking
# strings (e.g. "0.997") are not reinterpreted as numbers by Excel. ---
$priceCells = @("D2","D3","D4","D5","D6","D8","D10","D11","D12","D13","D15","D16","D17","D18","D19","D22","D23","D25","D27","D28","D31","D33","D36","D43","D46","D47","D48","D49","D50")
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "27.769.66"
$ws.Range("D3").Value = "1.630.75"
$ws.Range("D4").Value = "0.997"
$ws.Range("D5").Value = "211.67"
$ws.Range("D6").Value = "0.522"
$ws.Range("D8").Value = "23.24"
$ws.Range("D10").Value = "0.0614"
$ws.Range("D11").Value = "0.0881"
$ws.Range("D12").Value = "1.862.33"
$ws.Range("D13").Value = "1.631.19"
$ws.Range("D15").Value = "0.564"
$ws.Range("D16").Value = "65.31"
$ws.Range("D17").Value = "27.807.99"
$ws.Range("D18").Value = "230.43"
$ws.Range("D19").Value = "0.0₃0723"
$ws.Range("D22").Value = "10.37"
$ws.Range("D23").Value = "4.36"
$ws.Range("D25").Value = "154.53"
$ws.Range("D27").Value = "15.64"
$ws.Range("D28").Value = "0.111"
$ws.Range("D31").Value = "0.0482"
$ws.Range("D33").Value = "1.406.94"
$ws.Range("D36").Value = "1.01"
$ws.Range("D43").Value = "66.62"
$ws.Range("D46").Value = "2.19"
$ws.Range("D47").Value = "1.773.32"
$ws.Range("D48").Value = "87.95"
$ws.Range("D49").Value = "0.0₆0104"
$ws.Range("D50").Value = "0.0999"

foreach ($addr in $priceCells) {
    $ws.Range($addr).Style = "Normal"
}

# --- Coin name / link / volume updates ---
$ws.Range("E2").Value = "  -0.62%  "
$ws.Range("E3").Value = "  -0.90%  "
$ws.Range("E4").Value = "  -0.35%  "
$ws.Range("E6").Value = "  -0.94%  "
$ws.Range("E8").Value = "  -0.89%  "
$ws.Range("E9").Value = "  -2.70%  "
$ws.Range("E10").Value = "  -0.13%  "
$ws.Range("E11").Value = "  +0.99%  "
$ws.Range("E12").Value = "  -0.86%  "
$ws.Range("E13").Value = "  -0.76%  "
$ws.Range("E14").Value = "  -0.66%  "
$ws.Range("E15").Value = "  -0.16%  "
$ws.Range("E16").Value = "  -0.36%  "
$ws.Range("E17").Value = "  -0.48%  "
$ws.Range("E18").Value = "  -0.39%  "
$ws.Range("E19").Value = "  -0.18%  "
$ws.Range("E20").Value = "  -1.72%  "
$ws.Range("E21").Value = "  -0.31%  "
$ws.Range("B22").Value = "Avalanche"
$ws.Range("C22").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("E22").Value = "  -2.73%  "
$ws.Range("B23").Value = "Uniswap"
$ws.Range("C23").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("E23").Value = "  -0.77%  "
$ws.Range("E24").Value = "  -4.04%  "
$ws.Range("E25").Value = "  +1.43%  "
$ws.Range("E26").Value = "  +0.44%  "
$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("E27").Value = "  -0.70%  "
$ws.Range("B28").Value = "Stellar"
$ws.Range("C28").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("E28").Value = "  -0.98%  "
$ws.Range("E29").Value = "  -0.27%  "
$ws.Range("E30").Value = "  -1.23%  "
$ws.Range("E31").Value = "  -0.70%  "
$ws.Range("E32").Value = "  +2.51%  "
$ws.Range("E33").Value = "  -2.44%  "
$ws.Range("E34").Value = "  +0.31%  "
$ws.Range("E35").Value = "  +0.06%  "
$ws.Range("E36").Value = "  +7.90%  "
$ws.Range("E37").Value = "  +1.13%  "
$ws.Range("E38").Value = "  +0.64%  "
$ws.Range("E39").Value = "  +0.06%  "
$ws.Range("E40").Value = "  -2.49%  "
$ws.Range("E41").Value = "  -0.16%  "
$ws.Range("E42").Value = "  -0.22%  "
$ws.Range("E43").Value = "  -3.64%  "
$ws.Range("E45").Value = "  +0.18%  "
$ws.Range("E46").Value = "  -0.88%  "
$ws.Range("E47").Value = "  -0.84%  "
$ws.Range("E48").Value = "  -1.21%  "
$ws.Range("E49").Value = "  -1.40%  "
$ws.Range("E50").Value = "  -0.94%  "
$ws.Range("E51").Value = "  -0.23%  "
